# Weekly crypto price/volume refresh (GitHub Actions scrape), Tue Jun 11 2024 05:26:25 UTC
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The new Price values below are plain decimal numbers (single ".").
# Force the cell format to Text first so they stay literal strings, just
# like every other cell in this column, instead of becoming numeric values.
$textRows = @(5,6,11,15,19,20,21,22,24,27,29,30,31,33,34,36,38,40,42,44,45,46,47,48,49)
foreach ($r in $textRows) {
    $ws.Range("D$r").NumberFormat = "@"
}

$ws.Range("D2").Value = '67.956.71'
$ws.Range("E2").Value = '  -2.37%  '
$ws.Range("D3").Value = '3.564.22'
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '619.49'
$ws.Range("E5").Value = '  -7.11%  '
$ws.Range("D6").Value = '154.19'
$ws.Range("E6").Value = '  -3.59%  '
$ws.Range("D7").Value = '3.559.86'
$ws.Range("E7").Value = '  -3.37%  '
$ws.Range("E8").Value = '  -0.11%  '
$ws.Range("E9").Value = '  -2.16%  '
$ws.Range("E10").Value = '  -3.13%  '
$ws.Range("D11").Value = '6.91'
$ws.Range("E11").Value = '  -3.01%  '
$ws.Range("E12").Value = '  -1.80%  '
$ws.Range("E13").Value = '  -3.39%  '
$ws.Range("D14").Value = '4.166.54'
$ws.Range("E14").Value = '  -3.36%  '
$ws.Range("D15").Value = '32.07'
$ws.Range("E15").Value = '  -2.23%  '
$ws.Range("D16").Value = '3.566.06'
$ws.Range("E16").Value = '  -4.03%  '
$ws.Range("D17").Value = '67.993.42'
$ws.Range("E17").Value = '  -2.31%  '
$ws.Range("E18").Value = '  -1.11%  '
$ws.Range("D19").Value = '6.43'
$ws.Range("D20").Value = '15.65'
$ws.Range("E20").Value = '  -2.95%  '
$ws.Range("D21").Value = '456.35'
$ws.Range("E21").Value = '  -2.63%  '
$ws.Range("D22").Value = '9.67'
$ws.Range("E22").Value = '  -1.64%  '
$ws.Range("E23").Value = '  +0.10%  '
$ws.Range("D24").Value = '77.74'
$ws.Range("E24").Value = '  -2.65%  '
$ws.Range("D25").Value = '3.707.22'
$ws.Range("E25").Value = '  -3.32%  '
$ws.Range("E26").Value = '  +0.17%  '
$ws.Range("D27").Value = '10.68'
$ws.Range("E27").Value = '  -2.16%  '
$ws.Range("E28").Value = '  -7.88%  '
$ws.Range("D29").Value = '8.36'
$ws.Range("E29").Value = '  -7.58%  '
$ws.Range("D30").Value = '2.56'
$ws.Range("E30").Value = '  -4.25%  '
$ws.Range("D31").Value = '1.63'
$ws.Range("E31").Value = '  -3.73%  '
$ws.Range("E32").Value = '  -0.47%  '
$ws.Range("D33").Value = '25.99'
$ws.Range("E33").Value = '  -2.80%  '
$ws.Range("D34").Value = '1.91'
$ws.Range("E34").Value = '  -4.34%  '
$ws.Range("E35").Value = '  -4.08%  '
$ws.Range("B36").Value = 'NEARProtocol'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D36").Value = '6.20'
$ws.Range("E36").Value = '  -4.18%  '
$ws.Range("B37").Value = 'RenzoRestakedETH'
$ws.Range("C37").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D37").Value = '3.563.01'
$ws.Range("E37").Value = '  -3.16%  '
$ws.Range("D38").Value = '8.10'
$ws.Range("E38").Value = '  -3.98%  '
$ws.Range("E39").Value = '  -0.04%  '
$ws.Range("D40").Value = '178.15'
$ws.Range("E40").Value = '  -0.53%  '
$ws.Range("E41").Value = '  -0.09%  '
$ws.Range("D42").Value = '0.0888'
$ws.Range("E43").Value = '  -7.58%  '
$ws.Range("D44").Value = '2.10'
$ws.Range("E44").Value = '  -6.46%  '
$ws.Range("D45").Value = '0.895'
$ws.Range("E45").Value = '  -4.05%  '
$ws.Range("B46").Value = 'OKB'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D46").Value = '45.98'
$ws.Range("E46").Value = '  -2.20%  '
$ws.Range("B47").Value = 'InjectiveProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D47").Value = '28.60'
$ws.Range("E47").Value = '  +4.72%  '
$ws.Range("D48").Value = '2.60'
$ws.Range("E48").Value = '  -5.59%  '
$ws.Range("D49").Value = '7.74'
$ws.Range("E49").Value = '  -1.30%  '
$ws.Range("E50").Value = '  -5.70%  '
$ws.Range("E51").Value = '  -4.81%  '
